$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "310.97"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.38%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.48"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.42%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.110"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.61%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07781"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.92%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.374"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.218"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.77%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.886"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-7.77%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.751"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-9.98%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9226"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.38%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-5.13%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1911"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.69%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09198"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "3.91%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03429"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.92%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09678"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.21%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001368"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.91%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005991"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.39%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.557"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.62%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3372"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.76%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "4.54%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.79%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.47%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "5,592.05%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04362"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.01%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001208"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.68%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004251"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-8.72%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-63.79%"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-4.60%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05110"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.17%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007677"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.98%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009735"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.74%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1345"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.27%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.39%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009613"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "8.49%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006664"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.36%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.65%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.001201"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-0.68%"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002937"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-2.72%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.65%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.65%"
